$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value (45180 = 2023-09-11) for every
# data row (rows 2-484). The update bumps it by one day to 45181 (2023-09-12)
# for all rows.
$ws.Range("C2:C484").Value = 45181
